# Applies the "Add data for 2022-08-26" update:
#  - Renames the sheet / header label from "...August 17" to "...August 18"
#  - Updates the running-total counts in the "August 2022" column (B) and
#    several historical month columns with newly-arrived/corrected records.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "Through 2022-08-18"

# Update the column header text (column B, row 1) that mirrors the sheet name
$ws.Range("B1").Value = "August 2022 (through August 18)"

# Helper to add (or set) a cell value
function Set-CellValue($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Row 2 - Austin: new record in column AP, AO unchanged
Set-CellValue "AP2" 1

# Row 3 - Grand Crossing: R3 3 -> 4
Set-CellValue "R3" 4

# Row 4 - North Lawndale: B4 1 -> 2
Set-CellValue "B4" 2

# Row 5 - Garfield Park: Z5 2 -> 3
Set-CellValue "Z5" 3

# Row 7 - Humboldt Park: new record in column AH
Set-CellValue "AH7" 1

# Row 9 - Chatham: R9 3 -> 4
Set-CellValue "R9" 4

# Row 15 - Auburn Gresham: J15 1 -> 2
Set-CellValue "J15" 2

# Row 16 - Logan Square: R16 1 -> 2
Set-CellValue "R16" 2

# Row 21 - Wicker Park: new record in column B (August 2022)
Set-CellValue "B21" 1

# Row 38 - South Deering: new record in column AH
Set-CellValue "AH38" 1

# Row 39 - South Chicago: AP39 3 -> 4
Set-CellValue "AP39" 4

# Row 50 - Gage Park: AP50 1 -> 2
Set-CellValue "AP50" 2

# Row 53 - Bridgeport: new record in column B (August 2022)
Set-CellValue "B53" 1

# Row 57 - Albany Park: R57 4 -> 5
Set-CellValue "R57" 5

# Row 61 - Avondale: J61 2 -> 3
Set-CellValue "J61" 3

# Row 66 - Chicago Lawn: new record in column B (August 2022)
Set-CellValue "B66" 1

# Row 74 - Hyde Park: new record in column B (August 2022)
Set-CellValue "B74" 1
